# Update the "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages build at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 3129   # was 3120
$ws1.Range("F5").Value  = 6838   # was 6835
$ws1.Range("F6").Value  = 1899   # was 1895
$ws1.Range("F7").Value  = 8      # was 5
$ws1.Range("F8").Value  = 63     # was 61
$ws1.Range("F12").Value = 15     # was 12
$ws1.Range("F14").Value = 163    # was 161

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 3129   # was 3120
$ws4.Range("F6").Value  = 6838   # was 6835
$ws4.Range("F7").Value  = 1899   # was 1895
$ws4.Range("F8").Value  = 8      # was 5
$ws4.Range("F9").Value  = 63     # was 61
$ws4.Range("F13").Value = 15     # was 12
$ws4.Range("F15").Value = 163    # was 161
